# Update "想去人数" (interest count) figures in column F for the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets, refreshed from the
# upstream data source for the gh-pages output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 18533
$ws1.Range("F9").Value  = 1076
$ws1.Range("F10").Value = 6940
$ws1.Range("F11").Value = 414
$ws1.Range("F15").Value = 121
$ws1.Range("F17").Value = 225
$ws1.Range("F22").Value = 664
$ws1.Range("F26").Value = 286
$ws1.Range("F27").Value = 1013
$ws1.Range("F31").Value = 545
$ws1.Range("F33").Value = 66
$ws1.Range("F36").Value = 12191
$ws1.Range("F40").Value = 219
$ws1.Range("F41").Value = 300
$ws1.Range("F43").Value = 303

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 18533
$ws4.Range("F9").Value  = 1076
$ws4.Range("F10").Value = 6941
$ws4.Range("F11").Value = 414
$ws4.Range("F15").Value = 121
$ws4.Range("F17").Value = 225
$ws4.Range("F22").Value = 664
$ws4.Range("F26").Value = 286
$ws4.Range("F27").Value = 1013
$ws4.Range("F31").Value = 545
$ws4.Range("F35").Value = 66
$ws4.Range("F38").Value = 12191
$ws4.Range("F42").Value = 219
$ws4.Range("F43").Value = 300
$ws4.Range("F45").Value = 303
